$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new rows right after the header (row 1), pushing the
# existing data rows down, and fill them with the new translation
# test entries.
$ws.Range("A2:A4").EntireRow.Insert()

$ws.Range("A2").Value = "TEST TWO"
$ws.Range("A3").Value = "TEST ONE"
$ws.Range("A4").Value = "TEST THREE"

# The inserted rows should carry no special formatting (the header's
# blue-fill style must not bleed into them).
$ws.Range("A2:A4").ClearFormats()
